# Auto-generated Excel COM-interop script to apply scheduled Sheets update
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 16211.714
$ws.Range("I2").Value = 2246.6667
$ws.Range("J2").Value = 100002
$ws.Range("K2").Value = 2246.6667
$ws.Range("L2").Value = 100002
$ws.Range("M2").Value = -2133.6667
$ws.Range("N2").Value = -100228
$ws.Range("H9").Value = 1930.5
$ws.Range("I9").Value = 238.06667
$ws.Range("J9").Value = 7007.8
$ws.Range("K9").Value = 238.06667
$ws.Range("L9").Value = 7007.8
$ws.Range("M9").Value = -69.06666999999999
$ws.Range("N9").Value = -7345.8
$ws.Range("H33").Value = 11193.429
$ws.Range("I33").Value = 11193.429
$ws.Range("K33").Value = 11193.429
$ws.Range("M33").Value = -10964.429
$ws.Range("H40").Value = 2832.75
$ws.Range("I40").Value = 3220.5
$ws.Range("J40").Value = 2445
$ws.Range("K40").Value = 3220.5
$ws.Range("L40").Value = 2445
$ws.Range("M40").Value = -3045.5
$ws.Range("N40").Value = -2795
$ws.Range("H42").Value = 367
$ws.Range("I42").Value = 287.6
$ws.Range("J42").Value = 466.25
$ws.Range("K42").Value = 862.8000000000001
$ws.Range("L42").Value = 1398.75
$ws.Range("M42").Value = -632.8000000000001
$ws.Range("N42").Value = -1858.75
$ws.Range("H49").Value = 70
$ws.Range("I49").Value = 70
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 210
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -74
$ws.Range("N49").Value = $null   # was -509
$ws.Range("H51").Value = 21573.412
$ws.Range("J51").Value = 26363.545
$ws.Range("L51").Value = 26363.545
$ws.Range("N51").Value = -27331.545
$ws.Range("H101").Value = 1230.3334
$ws.Range("J101").Value = 493
$ws.Range("L101").Value = 1479
$ws.Range("N101").Value = -4723
$ws.Range("H125").Value = 2787
$ws.Range("I125").Value = 2787
$ws.Range("K125").Value = 25083
$ws.Range("M125").Value = -22623
$ws.Range("H135").Value = 1916.0571
$ws.Range("I135").Value = 885.4
$ws.Range("J135").Value = 8100
$ws.Range("K135").Value = 7968.599999999999
$ws.Range("L135").Value = 72900
$ws.Range("M135").Value = -5433.599999999999
$ws.Range("N135").Value = -77970
$ws.Range("H141").Value = 982
$ws.Range("I141").Value = 982
$ws.Range("K141").Value = 2946
$ws.Range("M141").Value = 2234

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 337.5
$ws.Range("I4").Value = 337.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 337.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -221.5
$ws.Range("N4").Value = $null   # was -1232
$ws.Range("H32").Value = 19914.389
$ws.Range("I32").Value = 19724.094
$ws.Range("K32").Value = 19724.094
$ws.Range("M32").Value = -19437.094
$ws.Range("H45").Value = 8172.8667
$ws.Range("I45").Value = 9532.75
$ws.Range("K45").Value = 9532.75
$ws.Range("M45").Value = -9155.75
$ws.Range("H58").Value = 24793.5
$ws.Range("J58").Value = 24793.5
$ws.Range("L58").Value = 24793.5
$ws.Range("N58").Value = -25653.5
$ws.Range("H74").Value = 374336.56
$ws.Range("J74").Value = 9480
$ws.Range("L74").Value = 9480
$ws.Range("N74").Value = -11228
$ws.Range("H77").Value = 374336.56
$ws.Range("J77").Value = 9480
$ws.Range("L77").Value = 47400
$ws.Range("N77").Value = -56136
$ws.Range("H132").Value = 8849.6
$ws.Range("I132").Value = 5642.4287
$ws.Range("K132").Value = 16927.2861
$ws.Range("M132").Value = -14397.2861

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1242.1
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1346
$ws.Range("H25").Value = 4232
$ws.Range("I25").Value = 196.8
$ws.Range("J25").Value = 6066.1816
$ws.Range("K25").Value = 196.8
$ws.Range("L25").Value = 6066.1816
$ws.Range("M25").Value = 38.19999999999999
$ws.Range("N25").Value = -6536.1816
$ws.Range("H37").Value = 2013
$ws.Range("I37").Value = 1026
$ws.Range("J37").Value = 3000
$ws.Range("K37").Value = 1026
$ws.Range("L37").Value = 3000
$ws.Range("M37").Value = -889
$ws.Range("N37").Value = -3274
$ws.Range("H76").Value = 9569.75
$ws.Range("J76").Value = 9569.75
$ws.Range("L76").Value = 9569.75
$ws.Range("N76").Value = -10199.75
$ws.Range("H79").Value = 9569.75
$ws.Range("J79").Value = 9569.75
$ws.Range("L79").Value = 9569.75
$ws.Range("N79").Value = -11753.75
$ws.Range("H99").Value = 2814.0588
$ws.Range("I99").Value = 2946
$ws.Range("J99").Value = 2497.4
$ws.Range("K99").Value = 2946
$ws.Range("L99").Value = 2497.4
$ws.Range("M99").Value = -1448
$ws.Range("N99").Value = -5493.4
$ws.Range("H126").Value = 49997.5
$ws.Range("J126").Value = 49997.5
$ws.Range("L126").Value = 49997.5
$ws.Range("N126").Value = -59877.5
$ws.Range("H134").Value = 4303.108
$ws.Range("I134").Value = 2069.7273
$ws.Range("K134").Value = 6209.1819
$ws.Range("M134").Value = -3674.1819

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29416366
$ws.Range("J31").Value = 7794.8887
$ws.Range("L31").Value = 7794.8887
$ws.Range("N31").Value = -8384.8887
$ws.Range("H34").Value = 29416366
$ws.Range("J34").Value = 7794.8887
$ws.Range("L34").Value = 7794.8887
$ws.Range("N34").Value = -8198.8887
$ws.Range("H58").Value = 4640.8438
$ws.Range("I58").Value = 3092.8572
$ws.Range("J58").Value = 7596.091
$ws.Range("K58").Value = 3092.8572
$ws.Range("L58").Value = 7596.091
$ws.Range("M58").Value = -2889.8572
$ws.Range("N58").Value = -8002.091
$ws.Range("H95").Value = 16570.428
$ws.Range("J95").Value = 16570.428
$ws.Range("L95").Value = 16570.428
$ws.Range("N95").Value = -22062.428
$ws.Range("H136").Value = 4640.8438
$ws.Range("I136").Value = 3092.8572
$ws.Range("J136").Value = 7596.091
$ws.Range("K136").Value = 9278.571599999999
$ws.Range("L136").Value = 22788.273
$ws.Range("M136").Value = -6728.571599999999
$ws.Range("N136").Value = -27888.273

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 670.3
$ws.Range("I132").Value = 670.3
$ws.Range("K132").Value = 6032.7
$ws.Range("M132").Value = -3502.7

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 23842.883
$ws.Range("I24").Value = 33285.715
$ws.Range("J24").Value = 17232.9
$ws.Range("K24").Value = 33285.715
$ws.Range("L24").Value = 17232.9
$ws.Range("M24").Value = -33112.715
$ws.Range("N24").Value = -17578.9
$ws.Range("H27").Value = 2002.75
$ws.Range("J27").Value = 2503.6667
$ws.Range("L27").Value = 2503.6667
$ws.Range("N27").Value = -2835.6667
$ws.Range("H92").Value = 20211.111
$ws.Range("J92").Value = 21487.5
$ws.Range("L92").Value = 21487.5
$ws.Range("N92").Value = -25231.5
$ws.Range("H134").Value = 137500
$ws.Range("J134").Value = 137500
$ws.Range("L134").Value = 412500
$ws.Range("N134").Value = -417570

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4505.891
$ws.Range("I46").Value = 1399.25
$ws.Range("J46").Value = 4801.7617
$ws.Range("K46").Value = 1399.25
$ws.Range("L46").Value = 4801.7617
$ws.Range("M46").Value = -1211.25
$ws.Range("N46").Value = -5177.7617
$ws.Range("H61").Value = 1203.2778
$ws.Range("I61").Value = 1086.1818
$ws.Range("K61").Value = 1086.1818
$ws.Range("M61").Value = -884.1818000000001
$ws.Range("H68").Value = 16479
$ws.Range("I68").Value = 3131.6667
$ws.Range("J68").Value = 36500
$ws.Range("K68").Value = 3131.6667
$ws.Range("L68").Value = 36500
$ws.Range("M68").Value = -2382.6667
$ws.Range("N68").Value = -37998
$ws.Range("H71").Value = 16479
$ws.Range("I71").Value = 3131.6667
$ws.Range("J71").Value = 36500
$ws.Range("K71").Value = 15658.3335
$ws.Range("L71").Value = 182500
$ws.Range("M71").Value = -11914.3335
$ws.Range("N71").Value = -189988
$ws.Range("H113").Value = 1203.2778
$ws.Range("I113").Value = 1086.1818
$ws.Range("K113").Value = 1086.1818
$ws.Range("M113").Value = 1083.8182

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null   # was -59799.5
